$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country name cells where ranking order changed (adjacent row swaps)
$ws.Range("A36").Value = "Ucrania"
$ws.Range("A37").Value = "Oman"
$ws.Range("A43").Value = "Rumania"
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("A68").Value = "Nepal"
$ws.Range("A69").Value = "Costa Rica"
$ws.Range("A78").Value = "Estado de Palestina"
$ws.Range("A79").Value = "Dinamarca"
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# Update statistic columns (B..H) with the latest data
$ws.Range("B14").Value = 328844
$ws.Range("C14").Value = 2132
$ws.Range("D14").Value = 286642
$ws.Range("E14").Value = 23586
$ws.Range("G14").Value = 189
$ws.Range("H14").Value = 18616
$ws.Range("B36").Value = 81957
$ws.Range("C36").Value = 1008
$ws.Range("D36").Value = 44359
$ws.Range("E36").Value = 35676
$ws.Range("G36").Value = 25
$ws.Range("H36").Value = 1922
$ws.Range("B37").Value = 81787
$ws.Range("C37").Value = 207
$ws.Range("D37").Value = 76124
$ws.Range("E37").Value = 5142
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 521
$ws.Range("B41").Value = 72400
$ws.Range("C41").Value = 687
$ws.Range("D41").Value = 64028
$ws.Range("E41").Value = 7890
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 482
$ws.Range("B43").Value = 62547
$ws.Range("C43").Value = 779
$ws.Range("D43").Value = 30311
$ws.Range("E43").Value = 29507
$ws.Range("G43").Value = 29
$ws.Range("H43").Value = 2729
$ws.Range("B44").Value = 62525
$ws.Range("D44").Value = 56568
$ws.Range("E44").Value = 5600
$ws.Range("H44").Value = 357
$ws.Range("B57").Value = 37162
$ws.Range("C57").Value = 108
$ws.Range("D57").Value = 26228
$ws.Range("E57").Value = 9606
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 1328
$ws.Range("B58").Value = 36708
$ws.Range("C58").Value = 105
$ws.Range("E58").Value = 2421
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1987
$ws.Range("B68").Value = 23310
$ws.Range("C68").Value = 338
$ws.Range("D68").Value = 16493
$ws.Range("E68").Value = 6738
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 79
$ws.Range("B69").Value = 23286
$ws.Range("D69").Value = 7730
$ws.Range("E69").Value = 15321
$ws.Range("H69").Value = 235
$ws.Range("B78").Value = 14510
$ws.Range("C78").Value = 302
$ws.Range("D78").Value = 8045
$ws.Range("E78").Value = 6365
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 100
$ws.Range("B79").Value = 14442
$ws.Range("D79").Value = 12840
$ws.Range("E79").Value = 985
$ws.Range("H79").Value = 617
$ws.Range("B82").Value = 13202
$ws.Range("C82").Value = 116
$ws.Range("D82").Value = 11011
$ws.Range("E82").Value = 2040
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 151
$ws.Range("B85").Value = 11312
$ws.Range("C85").Value = 137
$ws.Range("D85").Value = 7390
$ws.Range("E85").Value = 3686
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 236
$ws.Range("E95").Value = 288
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 333
$ws.Range("B129").Value = 2255
$ws.Range("C129").Value = 6
$ws.Range("D129").Value = 1960
$ws.Range("E129").Value = 167
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 128
$ws.Range("B133").Value = 1962
$ws.Range("C133").Value = 4
$ws.Range("D133").Value = 1838
$ws.Range("B158").Value = 847
$ws.Range("C158").Value = 6
$ws.Range("E158").Value = 435
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 13:51"
